# "ui and price fix"
# Adds two new pricing sections ("Terasse" and "Gulve") below the existing
# "Tag" (roof) section, and moves the active selection to C6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 33 - "Terasse" (data row with factor formulas, 20/60 multipliers)
# ---------------------------------------------------------------------
$ws.Range("A33").Value = "Terasse"
$ws.Range("B33").Value = 10000
$ws.Range("C33").Value = 3000
$ws.Range("D33").Value = 0.6
$ws.Range("E33").Value = 1.4

$ws.Range("F33").Formula = '=($C33*20*D33)+$B33'
$ws.Range("G33").Formula = '=($C33*20*E33)+$B33'
$ws.Range("H33").Formula = '=($C33*60*D33)+$B33'
$ws.Range("I33").Formula = '=($C33*60*E33)+$B33'
$ws.Range("F33:I33").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Row 34 - "faktor ved hævet"
# ---------------------------------------------------------------------
$ws.Range("A34").Value = "faktor ved hævet"
$ws.Range("C34").Value = 1.5

# ---------------------------------------------------------------------
# Row 35 - "faktor  ved værn"
# ---------------------------------------------------------------------
$ws.Range("A35").Value = "faktor  ved værn"
$ws.Range("C35").Value = 1.2

# ---------------------------------------------------------------------
# Row 36 - "tilvalg trappe"
# ---------------------------------------------------------------------
$ws.Range("A36").Value = "tilvalg trappe"
$ws.Range("B36").Value = 20000

# ---------------------------------------------------------------------
# Row 39 - "Gulve" (data row with factor formulas, 80/200 multipliers)
# ---------------------------------------------------------------------
$ws.Range("A39").Value = "Gulve"
$ws.Range("B39").Value = 10000
$ws.Range("C39").Value = 1000
$ws.Range("D39").Value = 0.4
$ws.Range("E39").Value = 2

$ws.Range("F39").Formula = '=($C39*80*D39)+$B39'
$ws.Range("G39").Formula = '=($C39*80*E39)+$B39'
$ws.Range("H39").Formula = '=($C39*200*D39)+$B39'
$ws.Range("I39").Formula = '=($C39*200*E39)+$B39'
$ws.Range("F39:I39").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Row 40 - "tillæg gulvvarme"
# ---------------------------------------------------------------------
$ws.Range("A40").Value = "tillæg gulvvarme"
$ws.Range("C40").Value = 500
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 1

$ws.Range("F40").Formula = '=($C40*80*D40)+$B40'
$ws.Range("G40").Formula = '=($C40*80*E40)+$B40'
$ws.Range("H40").Formula = '=($C40*200*D40)+$B40'
$ws.Range("I40").Formula = '=($C40*200*E40)+$B40'
$ws.Range("F40:I40").NumberFormat = "#,##0"

# ---------------------------------------------------------------------
# Move the active selection, matching the saved UI state in the workbook.
# ---------------------------------------------------------------------
$ws.Range("C6").Select()
